# PriceProf.xlsx data-analysis update
#
# - Row-2 price values on the six "xx-wd"/"xx-we" demand-response sheets
#   (Winter-wd, Winter-we, Summer-wd, Summer-we, Autumn-wd, Autumn-we) are
#   overwritten with the freshly computed hourly price profile (same 24
#   values, B2:Y2, on every sheet - matching what Spring-wd/Spring-we
#   already contained).
# - Every sheet's selection is moved from the old "whole sheet" / stray
#   selection onto the header/price row (A2:XFD2).
# - The active tab moves from Spring-wd to Autumn-we (last sheet).

$wb = $excel.ActiveWorkbook

# The 24 new hourly values (columns B..Y) shared by every refreshed sheet.
$newValues = @(
    48.944624780929573,
    43.622658697783599,
    39.785869761391325,
    35.339140722213898,
    34.399946907223317,
    38.497734686852816,
    46.68209985724301,
    56.86108389365085,
    61.509879300107698,
    64.005432747320356,
    64.479988390425845,
    65.1025028684668,
    63.0149505850612,
    59.302777027107354,
    56.521452036326195,
    53.889464876308658,
    53.659291613963418,
    61.521851209167799,
    68.095062982422306,
    68.305684417863745,
    62.993558244454753,
    58.1847630441682,
    60.070880904817756,
    55.209300444730303
)

# Sheets whose B2:Y2 row needs the refreshed price values.
$sheetsToRefresh = @("Winter-wd", "Winter-we", "Summer-wd", "Summer-we", "Autumn-wd", "Autumn-we")

foreach ($sheetName in $sheetsToRefresh) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($col = 2; $col -le 25; $col++) {
        $ws.Cells.Item(2, $col).Value = $newValues[$col - 2]
    }
}

# Update the selection on every sheet to A2:XFD2 (the price row), in
# workbook tab order, finishing on Autumn-we so it ends up the active tab.
$sheetOrder = @("Winter-wd", "Winter-we", "Spring-wd", "Spring-we", "Summer-wd", "Summer-we", "Autumn-wd", "Autumn-we")

foreach ($sheetName in $sheetOrder) {
    if ($sheetName -eq "Spring-we") {
        # Spring-we keeps its existing selection untouched.
        continue
    }
    $ws = $wb.Worksheets.Item($sheetName)
    [void]$ws.Range("A2:XFD2").Select()
}
